# Generate Report for Handback
#
# This applies the "handback" update to the localization-status workbook:
#  - Overview + per-language "Status" cells flip from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - The zh-cn and de-de sheets get their "Latest Target File" / "Latest
#    Handback File" / "Latest Handback DateTime" columns populated (they were
#    previously blank / sentinel) with the md/xlf hand-back filenames and a
#    handback timestamp, plus a hyperlink on the new "Latest Target File"
#    cells (mirroring the existing hyperlink style already used in column A).
#  - A few columns are widened to fit the newly-populated, longer text.

$wb = $excel.ActiveWorkbook

$handedBackText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status text for both languages
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $handedBackText
$overview.Range("F2").Value = $handedBackText
$overview.Range("E3").Value = $handedBackText
$overview.Range("F3").Value = $handedBackText

# Widen the two status columns to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# "Status" column mirrors the Overview status text.
$zhcn.Range("C2").Value = $handedBackText
$zhcn.Range("C3").Value = $handedBackText

# Row 2 (6ec12533...): Latest Target File / Latest Handback File / Latest Handback DateTime
$zhcn.Range("I2").Value = "6ec12533-a576-4f89-80bc-5822835220eb.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/ccc709a08dc65a88649444152ba5743f7b9c88e0/e2e/6ec12533-a576-4f89-80bc-5822835220eb.md", "", "", "6ec12533-a576-4f89-80bc-5822835220eb.md") | Out-Null
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("J2").Value = "6ec12533-a576-4f89-80bc-5822835220eb.f01e300f6df647d7e015df1058973a5e880b3cc6.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-13 14:58:31"

# Row 3 (88dd2981...): Latest Target File / Latest Handback File / Latest Handback DateTime
$zhcn.Range("I3").Value = "88dd2981-8365-455e-8597-10f036935f47.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/ccc709a08dc65a88649444152ba5743f7b9c88e0/e2e/88dd2981-8365-455e-8597-10f036935f47.md", "", "", "88dd2981-8365-455e-8597-10f036935f47.md") | Out-Null
$zhcn.Range("I3").Font.Underline = $true
$zhcn.Range("I3").Font.Color = 15570276
$zhcn.Range("J3").Value = "88dd2981-8365-455e-8597-10f036935f47.362778c489d34f051454bbf5fc6ea00466131be8.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-13 14:58:31"

# Widen columns to fit the newly populated / longer values.
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# "Status" column mirrors the Overview status text.
$dede.Range("C2").Value = $handedBackText
$dede.Range("C3").Value = $handedBackText

# Row 2 (6ec12533...): Latest Target File / Latest Handback File / Latest Handback DateTime
$dede.Range("I2").Value = "6ec12533-a576-4f89-80bc-5822835220eb.md"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/ccc709a08dc65a88649444152ba5743f7b9c88e0/e2e/6ec12533-a576-4f89-80bc-5822835220eb.md", "", "", "6ec12533-a576-4f89-80bc-5822835220eb.md") | Out-Null
$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = 15570276
$dede.Range("J2").Value = "6ec12533-a576-4f89-80bc-5822835220eb.f01e300f6df647d7e015df1058973a5e880b3cc6.de-de.xlf"
$dede.Range("K2").Value = "2016-08-13 14:58:41"

# Row 3 (88dd2981...): Latest Target File / Latest Handback File / Latest Handback DateTime
$dede.Range("I3").Value = "88dd2981-8365-455e-8597-10f036935f47.md"
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/ccc709a08dc65a88649444152ba5743f7b9c88e0/e2e/88dd2981-8365-455e-8597-10f036935f47.md", "", "", "88dd2981-8365-455e-8597-10f036935f47.md") | Out-Null
$dede.Range("I3").Font.Underline = $true
$dede.Range("I3").Font.Color = 15570276
$dede.Range("J3").Value = "88dd2981-8365-455e-8597-10f036935f47.362778c489d34f051454bbf5fc6ea00466131be8.de-de.xlf"
$dede.Range("K3").Value = "2016-08-13 14:58:41"

# Widen columns to fit the newly populated / longer values.
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
